# Insert a new weekly price-report row at row 48 (pushing the existing
# rows 48..141 down to 49..142) and populate it with the new record's
# data, matching the "Fruta / hortaliza, semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 48..141 down to 49..142, leaving row 48 free for the new entry.
$ws.Rows.Item(48).Insert()

$ws.Range("A48").Value = 2
$ws.Range("B48").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C48").Value = "Coquimbo"
$ws.Range("D48").Value = 45246
$ws.Range("E48").Value = 4
$ws.Range("F48").Value = 100112030
$ws.Range("G48").Value = "Poroto granado"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 600
$ws.Range("K48").Value = 28000
$ws.Range("L48").Value = 30000
$ws.Range("M48").Value = 29000
$ws.Range("N48").Value = '$/caja 15 kilos'
$ws.Range("O48").Value = "Provincia de Limarí"
$ws.Range("P48").Value = 1933
$ws.Range("Q48").Value = 15
$ws.Range("R48").Value = "Hortaliza"
